# Weekly data refresh: insert 5 new price rows (Fecha = 44714) right above the
# existing "Barraganete" block that used to start at row 1106, shifting all
# subsequent rows down by 5 (old row N -> new row N+5). The sheet grows from
# A1:T1172 to A1:T1177.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 blank rows at 1106..1110, pushing the old data down.
$ws.Range("A1106:A1110").EntireRow.Insert()

# Columns that are constant across every data row in this sheet.
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$tipo        = "Fruta"
$productoId  = 100108
$producto    = "Tropicales y subtropicales"
$categoriaId = 100108006
$categoria   = "Plátano"
$unidad      = "`$/caja 20 kilos"
$origen      = "Ecuador"
$kgPorUnidad = 20
$fecha       = 44714

# New weekly rows: RowNum, Variedad, Calidad, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg
$newRows = @(
  @(1106, "Barraganete",      "Primera",        324,  17000, 18000, 17500, 875),
  @(1107, "Sin especificar",  "Maduro",         1000,  9000,  9500,  9280, 464),
  @(1108, "Sin especificar",  "Pintón",         2060,  8000, 10000,  8942, 447),
  @(1109, "Sin especificar",  "Primera Maduro", 1000, 11000, 11500, 11260, 563),
  @(1110, "Sin especificar",  "Primera Pintón", 3900, 10000, 12000, 11077, 554)
)

foreach ($r in $newRows) {
  $rowNum = $r[0]

  $ws.Cells.Item($rowNum, 1).Value  = $mercadoId
  $ws.Cells.Item($rowNum, 2).Value  = $mercado
  $ws.Cells.Item($rowNum, 3).Value  = $region
  $ws.Cells.Item($rowNum, 4).Value  = $fecha
  $ws.Cells.Item($rowNum, 5).Value  = $codreg
  $ws.Cells.Item($rowNum, 6).Value  = $tipo
  $ws.Cells.Item($rowNum, 7).Value  = $productoId
  $ws.Cells.Item($rowNum, 8).Value  = $producto
  $ws.Cells.Item($rowNum, 9).Value  = $categoriaId
  $ws.Cells.Item($rowNum, 10).Value = $categoria
  $ws.Cells.Item($rowNum, 11).Value = $r[1]
  $ws.Cells.Item($rowNum, 12).Value = $r[2]
  $ws.Cells.Item($rowNum, 13).Value = $r[3]
  $ws.Cells.Item($rowNum, 14).Value = $r[4]
  $ws.Cells.Item($rowNum, 15).Value = $r[5]
  $ws.Cells.Item($rowNum, 16).Value = $r[6]
  $ws.Cells.Item($rowNum, 17).Value = $unidad
  $ws.Cells.Item($rowNum, 18).Value = $origen
  $ws.Cells.Item($rowNum, 19).Value = $r[7]
  $ws.Cells.Item($rowNum, 20).Value = $kgPorUnidad
}
